# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.178.03"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = "'1.830.95"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = "'1.011"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = "'313.16"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').Value = "'0.4707"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = "'0.3665"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = "'0.07402"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = "'0.8809"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').Value = "'20.33"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = "'1.896.51"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.48%  '
$ws.Range('D13').Value = "'0.07657"
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Value = "'93.43"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = "'5.384"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = "'0.000008735"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('D19').Value = "'1.009"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').Value = "'27.588.61"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('D21').Value = "'14.63"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').Value = "'5.249"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.60%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = "'2.087.47"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.58%  '
$ws.Range('D25').Value = "'1.882"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('D26').Value = "'151.15"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').Value = "'2.131"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('D29').Value = "'5.186"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('D30').Value = "'116.64"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = "'0.08942"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('D32').Value = "'0.7458"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Value = "'1.166"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').Value = "'4.519"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').Value = "'2.943"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('D37').Value = "'2.559"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.30%  '
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').Value = "'0.01940"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = "'2.935"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('D43').Value = "'0.5266"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('D45').Value = "'8.379"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('D46').Value = "'0.4903"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').Value = "'10.39"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').Value = "'104.51"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('D50').Value = "'1.655"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').Value = "'0.06280"
$ws.Range('D51').Style = 'Normal'
